# Weekly update: insert the newest week's "Frutilla" (strawberry) price
# triple (Especial / Primera / Segunda) ahead of the existing rows.
#
# The sheet stores one row per quality grade per reporting week, with the
# most-recent week's rows living at the top of the data block (row 573).
# Adding a new week means inserting 3 fresh rows at row 573 and pushing the
# rest of the table down - everything else about each historical row
# (dates, prices, etc.) stays exactly as it was, just shifted down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current first data row of this block.
$ws.Rows("573:575").Insert()

# Shared/constant columns for every row in this sub-sheet.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

function Set-FrutillaRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New week (serial date 44610) - Especial / Primera / Segunda grades.
Set-FrutillaRow 573 44610 "Especial" 400 11500 12000 11750 1679
Set-FrutillaRow 574 44610 "Primera"  400 9500  10000 9750  1393
Set-FrutillaRow 575 44610 "Segunda"  300 7500  8000  7750  1107
